$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's worker-detail table (rows 16-26) previously listed 11 period
# rows across several different workers. The new "account statement" keeps
# only the data for worker EDWAR ALEXANDER MUÑOZ REYES (doc 1019060050),
# across two overdue periods (2506 and 2508).
#
# Delete rows 17-25 first: this removes the old middle rows while leaving the
# old row 26 (which carries the "last row" bottom-border styling) to shift up
# into row 17 - so both remaining data rows keep their correct borders/styles
# without us having to hand-roll style ids.
$ws.Range("17:25").Delete()

# Row 16 (first data row): CC / 1019060050 / EDWAR ALEXANDER MUÑOZ REYES / period 2506
$ws.Range("C16").Value = "1019060050"
$ws.Range("D16").Value = "EDWAR ALEXANDER MUÑOZ REYES"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 9490
$ws.Range("G16").Value = 1423500

# Row 17 (second / last data row, was old row 26 - same worker, new period 2508)
$ws.Range("C17").Value = "1019060050"
$ws.Range("D17").Value = "EDWAR ALEXANDER MUÑOZ REYES"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Summary fields above the table: total overdue value and worker count
$ws.Range("E11").Value = 66430
$ws.Range("C13").Value = 1

# Column D ("Nombre Trabajador") shrinks now that only one (shorter) name
# remains in the table instead of several longer ones.
$ws.Columns("D").ColumnWidth = 31.8
